$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$shopId = "da04f4ef-ffb0-11ea-ba65-065a10bcba76"

for ($i = 21; $i -le 33; $i++) {
    $row = $i - 6  # 21 -> 15, 22 -> 16, ..., 33 -> 27
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $shopId
    $ws.Cells.Item($row, 3).Formula = "=_xlfn.CONCAT(" + [char]34 + "INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('" + [char]34 + ", B" + $row + ", " + [char]34 + "'), LPAD(" + [char]34 + ", A" + $row + ", " + [char]34 + ", 7, '0'), 'dish'" + [char]34 + ", " + [char]34 + ");" + [char]34 + ")"
}

$ws.Range("C19").Select()
